# Auto-generated by diff analysis
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 41668252
$ws.Range("I32").Value = 333333340
$ws.Range("J32").Value = 1810.2858
$ws.Range("K32").Value = 333333340
$ws.Range("L32").Value = 1810.2858
$ws.Range("M32").Value = -333333014
$ws.Range("N32").Value = -2462.2858
$ws.Range("H58").Value = 6493865.5
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 6000
$ws.Range("N58").Value = -6300
$ws.Range("H64").Value = 2949.95
$ws.Range("I64").Value = 2836.303
$ws.Range("J64").Value = 3485.7144
$ws.Range("K64").Value = 2836.303
$ws.Range("L64").Value = 3485.7144
$ws.Range("M64").Value = -2588.303
$ws.Range("N64").Value = -3981.7144
$ws.Range("H67").Value = 2949.95
$ws.Range("I67").Value = 2836.303
$ws.Range("J67").Value = 3485.7144
$ws.Range("K67").Value = 2836.303
$ws.Range("L67").Value = 3485.7144
$ws.Range("M67").Value = -1978.303
$ws.Range("N67").Value = -5201.7144
$ws.Range("H103").Value = 2519.6
$ws.Range("I103").Value = 2519.6
$ws.Range("K103").Value = 7558.799999999999
$ws.Range("M103").Value = -6972.799999999999
$ws.Range("H129").Value = 828.6
$ws.Range("I129").Value = 491
$ws.Range("K129").Value = 1473
$ws.Range("M129").Value = 3527
$ws.Range("H132").Value = 22435.223
$ws.Range("I132").Value = 3451
$ws.Range("K132").Value = 10353
$ws.Range("M132").Value = -7823

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1601.463
$ws.Range("I74").Value = 1418.825
$ws.Range("J74").Value = 2123.2856
$ws.Range("K74").Value = 1418.825
$ws.Range("L74").Value = 2123.2856
$ws.Range("M74").Value = -544.825
$ws.Range("N74").Value = -3871.2856
$ws.Range("H77").Value = 1601.463
$ws.Range("I77").Value = 1418.825
$ws.Range("J77").Value = 2123.2856
$ws.Range("K77").Value = 7094.125
$ws.Range("L77").Value = 10616.428
$ws.Range("M77").Value = -2726.125
$ws.Range("N77").Value = -19352.428
$ws.Range("H112").Value = 42520
$ws.Range("J112").Value = 42520
$ws.Range("L112").Value = 42520
$ws.Range("N112").Value = -45474
$ws.Range("H132").Value = 9435891
$ws.Range("I132").Value = 13889982
$ws.Range("K132").Value = 41669946
$ws.Range("M132").Value = -41667416

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 33625
$ws.Range("J38").Value = 33625
$ws.Range("L38").Value = 33625
$ws.Range("N38").Value = -34457
$ws.Range("H86").Value = 3197.182
$ws.Range("I86").Value = 3038
$ws.Range("J86").Value = 3329.8333
$ws.Range("K86").Value = 3038
$ws.Range("L86").Value = 3329.8333
$ws.Range("M86").Value = -1915
$ws.Range("N86").Value = -5575.8333
$ws.Range("H89").Value = 3197.182
$ws.Range("I89").Value = 3038
$ws.Range("J89").Value = 3329.8333
$ws.Range("K89").Value = 15190
$ws.Range("L89").Value = 16649.1665
$ws.Range("M89").Value = -9574
$ws.Range("N89").Value = -27881.1665
$ws.Range("H112").Value = 47351.5
$ws.Range("J112").Value = 47351.5
$ws.Range("L112").Value = 47351.5
$ws.Range("N112").Value = -50305.5
$ws.Range("H134").Value = 2605.9106
$ws.Range("I134").Value = 1471.1364
$ws.Range("K134").Value = 4413.4092
$ws.Range("M134").Value = -1878.4092

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 333333820
$ws.Range("I32").Value = 333333820
$ws.Range("K32").Value = 333333820
$ws.Range("M32").Value = -333333504

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 17866.834
$ws.Range("I22").Value = 550.5
$ws.Range("K22").Value = 1651.5
$ws.Range("M22").Value = -1482.5
$ws.Range("H27").Value = 17866.834
$ws.Range("I27").Value = 550.5
$ws.Range("K27").Value = 1651.5
$ws.Range("M27").Value = -1549.5
$ws.Range("H34").Value = 2550
$ws.Range("J34").Value = 3300
$ws.Range("L34").Value = 9900
$ws.Range("N34").Value = -10068
$ws.Range("H39").Value = 2400
$ws.Range("J39").Value = 2400
$ws.Range("L39").Value = 7200
$ws.Range("N39").Value = -7788
$ws.Range("H46").Value = 1750
$ws.Range("J46").Value = 1750
$ws.Range("L46").Value = 5250
$ws.Range("N46").Value = -5432
$ws.Range("H74").Value = 10082.833
$ws.Range("J74").Value = 10082.833
$ws.Range("L74").Value = 30248.499
$ws.Range("N74").Value = -32370.499
$ws.Range("H77").Value = 10082.833
$ws.Range("J77").Value = 10082.833
$ws.Range("L77").Value = 90745.497
$ws.Range("N77").Value = -101353.497
$ws.Range("H86").Value = 996.2381
$ws.Range("I86").Value = 654.7778
$ws.Range("J86").Value = 1252.3334
$ws.Range("K86").Value = 1964.3334
$ws.Range("L86").Value = 3757.0002
$ws.Range("M86").Value = -778.3334
$ws.Range("N86").Value = -6129.0002
$ws.Range("H89").Value = 996.2381
$ws.Range("I89").Value = 654.7778
$ws.Range("J89").Value = 1252.3334
$ws.Range("K89").Value = 5893.000199999999
$ws.Range("L89").Value = 11271.0006
$ws.Range("M89").Value = 34.9998000000005
$ws.Range("N89").Value = -23127.0006
$ws.Range("H115").Value = 4050.6538
$ws.Range("I115").Value = 528
$ws.Range("J115").Value = 4191.56
$ws.Range("K115").Value = 1584
$ws.Range("L115").Value = 12574.68
$ws.Range("M115").Value = -409
$ws.Range("N115").Value = -14924.68

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1510.6
$ws.Range("I122").Value = 1584.2858
$ws.Range("J122").Value = 1338.6666
$ws.Range("K122").Value = 4752.857400000001
$ws.Range("L122").Value = 4015.9998
$ws.Range("M122").Value = -2302.857400000001
$ws.Range("N122").Value = -8915.9998
$ws.Range("H132").Value = 47621948
$ws.Range("I132").Value = 66668730
$ws.Range("K132").Value = 200006190
$ws.Range("M132").Value = -200003660
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2246.4211
$ws.Range("I16").Value = 1575
$ws.Range("J16").Value = 10079.667
$ws.Range("K16").Value = 1575
$ws.Range("L16").Value = 10079.667
$ws.Range("M16").Value = -1405
$ws.Range("N16").Value = -10419.667
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 3985
$ws.Range("I40").Value = 3985
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3985
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3849
$ws.Range("N40").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 3433.3914
$ws.Range("I132").Value = 2529.75
$ws.Range("J132").Value = 5498.857
$ws.Range("K132").Value = 7589.25
$ws.Range("L132").Value = 16496.571
$ws.Range("M132").Value = -5059.25
$ws.Range("N132").Value = -21556.571
$ws.Range("H136").Value = 2684
$ws.Range("I136").Value = 1699.7778
$ws.Range("J136").Value = 4898.5
$ws.Range("K136").Value = 5099.3334
$ws.Range("L136").Value = 14695.5
$ws.Range("M136").Value = -2549.3334
$ws.Range("N136").Value = -19795.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3778
$ws.Range("I62").Value = 10002
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 10002
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -9378
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3778
$ws.Range("I65").Value = 10002
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 50010
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -46890
$ws.Range("N65").Value = -21240
$ws.Range("H122").Value = 2859422.8
$ws.Range("I122").Value = 3573428
$ws.Range("J122").Value = 3402.5
$ws.Range("K122").Value = 10720284
$ws.Range("L122").Value = 10207.5
$ws.Range("M122").Value = -10717834
$ws.Range("N122").Value = -15107.5
$ws.Range("H126").Value = 3269347.8
$ws.Range("I126").Value = 4202614
$ws.Range("J126").Value = 2916.5
$ws.Range("K126").Value = 12607842
$ws.Range("L126").Value = 8749.5
$ws.Range("M126").Value = -12605372
$ws.Range("N126").Value = -13689.5
$ws.Range("H132").Value = 3108696.2
$ws.Range("I132").Value = 7248609.5
$ws.Range("J132").Value = 3761.375
$ws.Range("K132").Value = 21745828.5
$ws.Range("L132").Value = 11284.125
$ws.Range("M132").Value = -21743298.5
$ws.Range("N132").Value = -16344.125
